# Actualización automática de scrims_actualizado.xlsx (2025-07-27 15:07:37)
# Appends new scrim-match result rows to four worksheets, reproducing the
# rows added by the upstream automation (Google Forms -> Sheets style
# exporter). Each new row copies the existing column layout:
#   A:C  = Team A players (Brawlers)
#   D:F  = Team B players (Brawlers)
#   G    = Winner ("Equipo 1" / "Equipo 2" / "Empate")
#   H:M  = Discord handles of the six players
#   N    = submission timestamp (literal text, NOT a date)

$wb = $excel.ActiveWorkbook

function Add-ScrimRow {
    param(
        $ws,
        [int]$RowNum,
        [object[]]$Values
    )

    for ($c = 1; $c -le $Values.Count; $c++) {
        $ws.Cells.Item($RowNum, $c).Value = $Values[$c - 1]
    }

    # Thin border around every populated cell in the row (A:N), matching
    # the bordered/“card” look used by every existing results row.
    $rowRange = $ws.Range($ws.Cells.Item($RowNum, 1), $ws.Cells.Item($RowNum, 14))
    $rowRange.Borders.LineStyle = 1

    # Column G ("Ganador") is bold, like every other results row.
    $ws.Cells.Item($RowNum, 7).Font.Bold = $true
}

# ---------------------------------------------------------------------
# Sheet: Triple Dribble (rows 103-107)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Triple Dribble")

Add-ScrimRow $ws 103 @("HANK","CROW","BERRY","WILLOW","ASH","CHARLIE","Equipo 2","Tatsuki.💚","Yutapin","Shigemyon","FZ|Danshari","FZ|Toridesu","FZ|Mira","20250727T130654.000Z")
Add-ScrimRow $ws 104 @("HANK","CROW","BERRY","WILLOW","ASH","CHARLIE","Equipo 1","Tatsuki.💚","Yutapin","Shigemyon","FZ|Danshari","FZ|Toridesu","FZ|Mira","20250727T130427.000Z")
Add-ScrimRow $ws 105 @("HANK","CROW","BERRY","WILLOW","ASH","CHARLIE","Equipo 2","Tatsuki.💚","Yutapin","Shigemyon","FZ|Danshari","FZ|Toridesu","FZ|Mira","20250727T130233.000Z")
Add-ScrimRow $ws 106 @("ASH","WILLOW","JAE-YONG","MOE","LARRY & LAWRIE","BONNIE","Equipo 2","Tatsuki.💚","Yutapin","Shigemyon","FZ|Toridesu","FZ|Danshari","FZ|Mira","20250727T125530.000Z")
Add-ScrimRow $ws 107 @("ASH","WILLOW","JAE-YONG","MOE","LARRY & LAWRIE","BONNIE","Equipo 2","Tatsuki.💚","Yutapin","Shigemyon","FZ|Toridesu","FZ|Danshari","FZ|Mira","20250727T125231.000Z")

# ---------------------------------------------------------------------
# Sheet: Hot Potato (rows 109-112)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Hot Potato")

Add-ScrimRow $ws 109 @("WILLOW","BULL","LUMI","GRIFF","AMBER","BONNIE","Equipo 2","HMB|BosS","HMB|Symantec","IDarkLukii","SUP|Filippo神","SUP|Tomzy","Enraged 💔","20250727T130511.000Z")
Add-ScrimRow $ws 110 @("WILLOW","BULL","LUMI","GRIFF","AMBER","BONNIE","Equipo 1","HMB|BosS","HMB|Symantec","IDarkLukii","SUP|Filippo神","SUP|Tomzy","Enraged 💔","20250727T130220.000Z")
Add-ScrimRow $ws 111 @("AMBER","LILY","LOU","CARL","EMZ","BULL","Equipo 1","HMB|BosS","HMB|Symantec","IDarkLukii","Enraged 💔","SUP|Tomzy","SUP|Filippo神","20250727T125625.000Z")
Add-ScrimRow $ws 112 @("AMBER","LILY","LOU","CARL","EMZ","BULL","Equipo 1","HMB|BosS","HMB|Symantec","IDarkLukii","Enraged 💔","SUP|Tomzy","SUP|Filippo神","20250727T125340.000Z")

# ---------------------------------------------------------------------
# Sheet: Layer Cake (rows 108-111)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Layer Cake")

Add-ScrimRow $ws 108 @("JAE-YONG","BUSTER","LUMI","KIT","DOUG","EMZ","Equipo 1","HMB|BosS","HMB|Symantec","IDarkLukii","Enraged 💔","SUP|Tomzy","SUP|Filippo神","20250727T124631.000Z")
Add-ScrimRow $ws 109 @("JAE-YONG","BUSTER","LUMI","KIT","DOUG","EMZ","Equipo 1","HMB|BosS","HMB|Symantec","IDarkLukii","Enraged 💔","SUP|Tomzy","SUP|Filippo神","20250727T124425.000Z")
Add-ScrimRow $ws 110 @("LUMI","OLLIE","JAE-YONG","KIT","DOUG","SHADE","Equipo 1","IDarkLukii","HMB|Symantec","HMB|BosS","Enraged 💔","SUP|Filippo神","SUP|Tomzy","20250727T123810.000Z")
Add-ScrimRow $ws 111 @("LUMI","OLLIE","JAE-YONG","KIT","DOUG","SHADE","Equipo 1","IDarkLukii","HMB|Symantec","HMB|BosS","Enraged 💔","SUP|Filippo神","SUP|Tomzy","20250727T123603.000Z")

# ---------------------------------------------------------------------
# Sheet: Open Business (rows 132-133)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Open Business")

Add-ScrimRow $ws 132 @("FRANK","KIT","BEA","DRACO","WILLOW","MEG","Equipo 1","Shigemyon","Tatsuki.💚","Yutapin","FZ|Toridesu","FZ|Danshari","FZ|Mira","20250727T124411.000Z")
Add-ScrimRow $ws 133 @("FRANK","KIT","BEA","DRACO","WILLOW","MEG","Equipo 1","Shigemyon","Tatsuki.💚","Yutapin","FZ|Toridesu","FZ|Danshari","FZ|Mira","20250727T124144.000Z")

Write-Host "Scrim rows appended."
